# The "Sorting Results" sheet used to repeat the header row before every
# algorithm's data row (rows 1-12: header/BubbleSort/header/InsertionSort/...).
# Collapse that back down to a single header followed by one row per
# algorithm, and bump the "Data size" column from 3 to 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old (duplicated-header) table so no stale rows 8-12 survive.
$ws.Range("A1:C12").Clear()

$data = @(
    @("Sort Algorithm", "Data size", "Time Taken"),
    @("BubbleSort",     "'5", "'0"),
    @("InsertionSort",  "'5", "'0"),
    @("MergeSort",      "'5", "'0"),
    @("QuickSort",      "'5", "'0"),
    @("SelectionSort",  "'5", "'0"),
    @("ShellSort",      "'5", "'0")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
}

# Leading apostrophes above force "5"/"0" to be stored as text (matching
# the original file, where every cell -- including the numeric-looking
# ones -- is a shared string), but they also stamp a quote-prefix number
# format on the cells. Strip that back off so formatting stays untouched.
$ws.Range("A1:C7").ClearFormats()
